# Stellar Horizons SoP cards - correction
# 1) Bump the cached datetimeFigureOut fields (Date Placeholder) in the
#    slide master, every slide layout and the notes master forward one
#    day: 17/12/24 -> 18/12/24 (en-IT "D/M/YY") and 12/17/24 -> 12/18/24
#    (en-US "M/D/YY").
# 2) Resize/reposition the "Roll for resource production" note box on
#    slide 4 and correct its wording.

$p = $ppt.ActivePresentation

function Update-DateShape {
    param($shape)

    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    if ($tr.Length -le 0) { return }

    $hit = $tr.Find("12/17/24", 0)
    if ($hit) {
        $hit.Text = "12/18/24"
        return
    }

    $hit = $tr.Find("17/12/24", 0)
    if ($hit) {
        $hit.Text = "18/12/24"
        return
    }
}

# -- Slide master --
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# -- Every slide layout --
$layouts = $master.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# -- Notes master --
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    Update-DateShape $notesMaster.Shapes.Item($i)
}

# -- Slide 4: "Roll for resource production" note box --
$slide4 = $p.Slides.Item(4)
$box = $slide4.Shapes.Item(21)

if ($box.Name -ne "TextBox 58") {
    for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
        $cand = $slide4.Shapes.Item($i)
        if ($cand.Name -eq "TextBox 58") {
            $box = $cand
        }
    }
}

# Reposition / resize (points; values chosen so the Single-precision
# round trip lands exactly on the target EMU offsets/extents).
$box.Left = 362.40153506299214
$box.Top = 53.04515841023622
$box.Width = 178.5608215015748
$box.Height = 201.1453094905512

$tr = $box.TextFrame.TextRange

$r1 = $tr.Find("Roll for resource production", 0)
if ($r1) {
    $r1.Text = "Roll 1d10 for resource production"
}

$r2 = $tr.Find(". You always produce at least one.", 0)
if ($r2) {
    $r2.Text = ". You always produce at least one. +1 if modified roll <= world production value. +1 if mod roll = 0 or 1"
}
